# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# sheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

$updates = @{
    "F2" = 1895
    "F3" = 358
    "F4" = 1161
    "F5" = 1196
    "F7" = 5997
    "F8" = 100
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
